$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.454.94"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "1.909.66"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  +0.71%  "

$ws.Range("D5").Value = "'325.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "

$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("D7").Value = "'0.4812"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("D8").Value = "'0.4062"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.08157"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").Value = "'1.013"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'23.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.93%  "

$ws.Range("D12").Value = "1.910.62"
$ws.Range("E12").Value = "  -1.82%  "

$ws.Range("D13").Value = "'6.005"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.98%  "

$ws.Range("D14").Value = "'7.163"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").Value = "'90.30"
$ws.Range("D15").Style = "Normal"

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.06794"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.42%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").Value = "'0.00001037"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").Value = "'17.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").Value = "'1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").Value = "29.470.49"

$ws.Range("D22").Value = "'5.612"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").Value = "'11.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.38%  "

$ws.Range("D24").Value = "'2.182"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").Value = "2.166.23"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "'155.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "

$ws.Range("D27").Value = "'6.415"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.18%  "

$ws.Range("D28").Value = "'20.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("D29").Value = "'2.110"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "'120.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.09%  "

$ws.Range("D31").Value = "'1.025"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.27%  "

$ws.Range("D32").Value = "'0.09533"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "'5.534"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.76%  "

$ws.Range("D34").Value = "'3.559"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.46%  "

$ws.Range("D35").Value = "'1.388"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.64%  "

$ws.Range("D36").Value = "'0.02269"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").Value = "'0.06098"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").Value = "'1.178"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").Value = "'10.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.12%  "

$ws.Range("D40").Value = "'0.5954"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.48%  "

$ws.Range("D41").Value = "'7.992"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("D42").Value = "'0.1852"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("D43").Value = "'1.279"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").Value = "'2.370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.69%  "

$ws.Range("D45").Value = "'12.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.07%  "

$ws.Range("D46").Value = "'0.07596"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").Value = "'0.5569"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").Value = "'1.943"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("D49").Value = "'116.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.48%  "

$ws.Range("D50").Value = "'72.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.61%  "

$ws.Range("D51").Value = "'2.403"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
